$d = $word.ActiveDocument
$wdNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1. Title
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Find.Execute("AI's Ubiquitous Presence: Past, Present, and Future", $true, $false, $false, $false, $false, $true, 1, $false, "The Symphony of Life: A Journey Through Biology", 2)

# 2. Author name paragraph: "Mia Harrison" -> "Dr" + "." + " Sarah Peterson" (3 runs, same rPr)
$p2 = $d.Paragraphs.Item(2)
$p2full = $p2.Range
$p2sub = $d.Range($p2full.Start, $p2full.End - 1)
$nameXml = '<w:p ' + $wdNS + '><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>Dr</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve"> Sarah Peterson</w:t></w:r></w:p>'
$p2sub.InsertXML($nameXml)

# 3. Email paragraph (3 text-bearing runs change, periods stay)
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Find.Execute("mia", $true, $false, $false, $false, $false, $true, 1, $false, "biology", 2)
$p3b = $d.Paragraphs.Item(3)
$p3b.Range.Find.Execute("harrison@emailhost", $true, $false, $false, $false, $false, $true, 1, $false, "explorer123@eduverse", 2)
$p3c = $d.Paragraphs.Item(3)
$p3c.Range.Find.Execute("com", $true, $false, $false, $false, $false, $true, 1, $false, "org", 2)

# 4. Body paragraph: full restructure (text changes + 4 new runs + 2 runs consolidated to 1)
$p5 = $d.Paragraphs.Item(5)
$p5full = $p5.Range
$p5sub = $d.Range($p5full.Start, $p5full.End - 1)
$RPR = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr>'
$bodyXml = '<w:p ' + $wdNS + '>'
$bodyXml += '<w:r>' + $RPR + '<w:t>Biology is the study of life, exploring the remarkable intricacies of living organisms</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t xml:space="preserve"> It delves into the very essence of existence, questioning how life began, evolved, and sustains itself</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t xml:space="preserve"> As a science, biology unveils the secrets hidden within the vast tapestry of life, unraveling the mysteries that govern the behavior of organisms</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t xml:space="preserve"> Whether it be a tiny bacterium or a towering sequoia, biology seeks to understand the intricate mechanisms that orchestrate the symphony of life</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:br/></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:br/><w:t>The field of biology encompasses a broad spectrum of sub-disciplines, each specializing in a specific facet of life</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t xml:space="preserve"> Molecular biology explores the minute details of molecules and their role in biological processes</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t xml:space="preserve"> Cellular biology delves into the structure and function of cells, the fundamental building blocks of life</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t xml:space="preserve"> Evolution probes the transformative history of organisms, tracing the remarkable adaptations that have occurred over time</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t xml:space="preserve"> From the anatomy and physiology of organisms to their interactions with the environment, biology provides a comprehensive examination of the mechanisms that govern living systems</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:br/></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:br/><w:t>Biology is not merely an academic pursuit; it holds immense significance in our everyday lives</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t xml:space="preserve"> From the food we consume to the medicines we rely on, biology plays a crucial role in human well-being</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t xml:space="preserve"> Understanding the principles of biology empowers us to make informed decisions about our health, environment, and food choices</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t xml:space="preserve"> Whether we are studying the behavior of animals, unraveling the mysteries of human physiology, or developing innovative medical treatments, biology lies at the heart of our understanding of the living world and our place within it</w:t></w:r>'
$bodyXml += '<w:r>' + $RPR + '<w:t>.</w:t></w:r>'
$bodyXml += '</w:p>'
$p5sub.InsertXML($bodyXml)

# 5. Summary paragraph: 1:1 sentence replacements
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Find.Execute("Artificial Intelligence has woven itself into the fabric of our existence, serving as a tool of convenience, efficiency, and unparalleled potential", $true, $false, $false, $false, $false, $true, 1, $false, "Biology is the study of life, unravelling the intricate mechanisms that orchestrate the symphony of organisms", 2)
$p7b = $d.Paragraphs.Item(7)
$p7b.Range.Find.Execute(" Its journey from theoretical concept to ubiquitous reality is a testament to its transformative power", $true, $false, $false, $false, $false, $true, 1, $false, " With a vast tapestry of sub-disciplines, biology explores the complexities of living systems, from the molecular level to the behavior of entire ecosystems", 2)
$p7c = $d.Paragraphs.Item(7)
$p7c.Range.Find.Execute(" The synthesis of human ingenuity and AI's boundless capabilities holds the promise of addressing global ", $true, $false, $false, $false, $false, $true, 1, $false, " Its relevance extends beyond academia, impacting human well-being, ", 2)
$p7d = $d.Paragraphs.Item(7)
$p7d.Range.Find.Execute("challenges, redefining industries, and ushering in a new era of progress", $true, $false, $false, $false, $false, $true, 1, $false, "healthcare, and our understanding of the world around us", 2)
$p7e = $d.Paragraphs.Item(7)
$p7e.Range.Find.Execute(" As we embrace the future, AI stands poised to redefine the contours of human experience, opening doors to unimaginable possibilities", $true, $false, $false, $false, $false, $true, 1, $false, " By delving into the mysteries of life, biology empowers us to make informed decisions and fosters an appreciation for the delicate balance that sustains the planet we inhabit", 2)

# 6. Add a new empty paragraph at the very end of the document body
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML('<w:p ' + $wdNS + '/>')

Write-Host "Done"
